$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=FAPs, Target=ECs): update TPM-derived columns K..T
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08037899999999999
$ws.Range("N2").Value = 0.241137
$ws.Range("O2").Value = 0.01215705881223039
$ws.Range("P2").Value = 0.01215705881223039
$ws.Range("Q2").Value = 0.012673490895
$ws.Range("R2").Value = 0.114061418055
$ws.Range("S2").Value = 0.01215705881223039
$ws.Range("T2").Value = 0.01215705881223039

# Row 3 (Sending=FAPs, Target=FAPs): update recalculated specificity columns
$ws.Range("O3").Value = 0.2935162100923598
$ws.Range("P3").Value = 0.2935162100923598
$ws.Range("Q3").Value = 0.3059847841155556
$ws.Range("S3").Value = 0.2935162100923598
$ws.Range("T3").Value = 0.2935162100923598

# Row 4 (Sending=FAPs, Target=MuSCs): update TPM-derived columns K..T
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.590689999999999
$ws.Range("N4").Value = 13.77207
$ws.Range("O4").Value = 0.6943267310954097
$ws.Range("P4").Value = 0.6943267310954098
$ws.Range("Q4").Value = 0.7238217434499999
$ws.Range("R4").Value = 6.51439569105
$ws.Range("S4").Value = 0.6943267310954097
$ws.Range("T4").Value = 0.6943267310954098

# Row 5 (Sending=FAPs, Target=Resolving-Mac) no longer present in TPM data: remove it entirely
$ws.Rows.Item(5).Delete()
